# "Generate Report for Archive"
#
# The localization status of 0ce15954-64c9-4573-b9e5-bba00073ff90.md changed
# from "Ready for handoff" to "In Translation". The report is regenerated
# with rows grouped by status (In Translation rows first, sorted by file
# name, then Ready for handoff rows, then the fixed trailing rows), so
# 0ce15954... moves from row 4 up to row 2, pushing 78fed758... and
# a8aa8e86... down by one row each on every sheet. Row 5 (99319a3b...) and
# row 6 (.localization-config) are unaffected.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "0ce15954-64c9-4573-b9e5-bba00073ff90.md"
$ws1.Range("B2").Value = "In Translation"
$ws1.Range("C2").Value = "In Translation"

$ws1.Range("A3").Value = "78fed758-84e7-4089-a043-d6f60c88704b.md"
$ws1.Range("B3").Value = "In Translation"
$ws1.Range("C3").Value = "In Translation"

$ws1.Range("A4").Value = "a8aa8e86-e13d-403a-8a83-e5f5fb1a2894.md"
$ws1.Range("B4").Value = "In Translation"
$ws1.Range("C4").Value = "In Translation"

# Rows 5 and 6 are unchanged (99319a3b..., .localization-config).

$hls1 = $ws1.Hyperlinks
$hls1.Delete()
$hls1.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1be5524385e08f9374ae3ce9e24c96d6f00c8fa2/e2e/0ce15954-64c9-4573-b9e5-bba00073ff90.md", "", "", "0ce15954-64c9-4573-b9e5-bba00073ff90.md") | Out-Null
$hls1.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6906e8fdd8045ec51921315c59680c9f0e5ebab5/e2e/78fed758-84e7-4089-a043-d6f60c88704b.md", "", "", "78fed758-84e7-4089-a043-d6f60c88704b.md") | Out-Null
$hls1.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6906e8fdd8045ec51921315c59680c9f0e5ebab5/e2e/a8aa8e86-e13d-403a-8a83-e5f5fb1a2894.md", "", "", "a8aa8e86-e13d-403a-8a83-e5f5fb1a2894.md") | Out-Null
$hls1.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9c7cda88b83d9c32733cab49e364485c4000bf43/e2e/99319a3b-4a97-4338-b29b-ed4a462b98fc.md", "", "", "99319a3b-4a97-4338-b29b-ed4a462b98fc.md") | Out-Null
$hls1.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/9c7cda88b83d9c32733cab49e364485c4000bf43/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "0ce15954-64c9-4573-b9e5-bba00073ff90.md"
$ws2.Range("B2").Value = "In Translation"
$ws2.Range("C2").Value = "0ce15954-64c9-4573-b9e5-bba00073ff90.029dde67723be70b75b82033a51105d703a026b8.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-02-26 06:06:39"

$ws2.Range("A3").Value = "78fed758-84e7-4089-a043-d6f60c88704b.md"
$ws2.Range("B3").Value = "In Translation"
$ws2.Range("C3").Value = "78fed758-84e7-4089-a043-d6f60c88704b.6f5a7c6b436eb55169c871b262b35f9dde0a93e2.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-02-26 06:05:06"

$ws2.Range("A4").Value = "a8aa8e86-e13d-403a-8a83-e5f5fb1a2894.md"
$ws2.Range("B4").Value = "In Translation"
$ws2.Range("C4").Value = "a8aa8e86-e13d-403a-8a83-e5f5fb1a2894.01c7048035d0e34fda3ee283b56e83b8ec8c6c1a.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-02-26 06:05:06"

# Rows 5 and 6 are unchanged (99319a3b..., .localization-config).

$hls2 = $ws2.Hyperlinks
$hls2.Delete()
$hls2.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1be5524385e08f9374ae3ce9e24c96d6f00c8fa2/e2e/0ce15954-64c9-4573-b9e5-bba00073ff90.md", "", "", "0ce15954-64c9-4573-b9e5-bba00073ff90.md") | Out-Null
$hls2.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e862b984085792c2c20a027c35ee2ce6b30b62a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/0ce15954-64c9-4573-b9e5-bba00073ff90.029dde67723be70b75b82033a51105d703a026b8.zh-cn.xlf", "", "", "0ce15954-64c9-4573-b9e5-bba00073ff90.029dde67723be70b75b82033a51105d703a026b8.zh-cn.xlf") | Out-Null
$hls2.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6906e8fdd8045ec51921315c59680c9f0e5ebab5/e2e/78fed758-84e7-4089-a043-d6f60c88704b.md", "", "", "78fed758-84e7-4089-a043-d6f60c88704b.md") | Out-Null
$hls2.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d4fd8126926b28f057682196f72ddc8d7c5624d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/78fed758-84e7-4089-a043-d6f60c88704b.6f5a7c6b436eb55169c871b262b35f9dde0a93e2.zh-cn.xlf", "", "", "78fed758-84e7-4089-a043-d6f60c88704b.6f5a7c6b436eb55169c871b262b35f9dde0a93e2.zh-cn.xlf") | Out-Null
$hls2.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6906e8fdd8045ec51921315c59680c9f0e5ebab5/e2e/a8aa8e86-e13d-403a-8a83-e5f5fb1a2894.md", "", "", "a8aa8e86-e13d-403a-8a83-e5f5fb1a2894.md") | Out-Null
$hls2.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d4fd8126926b28f057682196f72ddc8d7c5624d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a8aa8e86-e13d-403a-8a83-e5f5fb1a2894.01c7048035d0e34fda3ee283b56e83b8ec8c6c1a.zh-cn.xlf", "", "", "a8aa8e86-e13d-403a-8a83-e5f5fb1a2894.01c7048035d0e34fda3ee283b56e83b8ec8c6c1a.zh-cn.xlf") | Out-Null
$hls2.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9c7cda88b83d9c32733cab49e364485c4000bf43/e2e/99319a3b-4a97-4338-b29b-ed4a462b98fc.md", "", "", "99319a3b-4a97-4338-b29b-ed4a462b98fc.md") | Out-Null
$hls2.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/68fa337acd41dcfd2e757a9892cbfbb0d2764e8e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/99319a3b-4a97-4338-b29b-ed4a462b98fc.91f6bf9d511a154dc72d665b05a013e920cdbb30.zh-cn.xlf", "", "", "99319a3b-4a97-4338-b29b-ed4a462b98fc.91f6bf9d511a154dc72d665b05a013e920cdbb30.zh-cn.xlf") | Out-Null
$hls2.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/9c7cda88b83d9c32733cab49e364485c4000bf43/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "0ce15954-64c9-4573-b9e5-bba00073ff90.md"
$ws3.Range("B2").Value = "In Translation"
$ws3.Range("C2").Value = "0ce15954-64c9-4573-b9e5-bba00073ff90.029dde67723be70b75b82033a51105d703a026b8.de-de.xlf"
$ws3.Range("D2").Value = "2016-02-26 06:06:51"

$ws3.Range("A3").Value = "78fed758-84e7-4089-a043-d6f60c88704b.md"
$ws3.Range("B3").Value = "In Translation"
$ws3.Range("C3").Value = "78fed758-84e7-4089-a043-d6f60c88704b.6f5a7c6b436eb55169c871b262b35f9dde0a93e2.de-de.xlf"
$ws3.Range("D3").Value = "2016-02-26 06:05:20"

$ws3.Range("A4").Value = "a8aa8e86-e13d-403a-8a83-e5f5fb1a2894.md"
$ws3.Range("B4").Value = "In Translation"
$ws3.Range("C4").Value = "a8aa8e86-e13d-403a-8a83-e5f5fb1a2894.01c7048035d0e34fda3ee283b56e83b8ec8c6c1a.de-de.xlf"
$ws3.Range("D4").Value = "2016-02-26 06:05:20"

# Rows 5 and 6 are unchanged (99319a3b..., .localization-config).

$hls3 = $ws3.Hyperlinks
$hls3.Delete()
$hls3.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1be5524385e08f9374ae3ce9e24c96d6f00c8fa2/e2e/0ce15954-64c9-4573-b9e5-bba00073ff90.md", "", "", "0ce15954-64c9-4573-b9e5-bba00073ff90.md") | Out-Null
$hls3.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c5f55c0858a5f01f79f04aabc50fadd29f4bd314/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/0ce15954-64c9-4573-b9e5-bba00073ff90.029dde67723be70b75b82033a51105d703a026b8.de-de.xlf", "", "", "0ce15954-64c9-4573-b9e5-bba00073ff90.029dde67723be70b75b82033a51105d703a026b8.de-de.xlf") | Out-Null
$hls3.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6906e8fdd8045ec51921315c59680c9f0e5ebab5/e2e/78fed758-84e7-4089-a043-d6f60c88704b.md", "", "", "78fed758-84e7-4089-a043-d6f60c88704b.md") | Out-Null
$hls3.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d4040b3e1f70ff1b3206d7be231a7cda26b7619b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/78fed758-84e7-4089-a043-d6f60c88704b.6f5a7c6b436eb55169c871b262b35f9dde0a93e2.de-de.xlf", "", "", "78fed758-84e7-4089-a043-d6f60c88704b.6f5a7c6b436eb55169c871b262b35f9dde0a93e2.de-de.xlf") | Out-Null
$hls3.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6906e8fdd8045ec51921315c59680c9f0e5ebab5/e2e/a8aa8e86-e13d-403a-8a83-e5f5fb1a2894.md", "", "", "a8aa8e86-e13d-403a-8a83-e5f5fb1a2894.md") | Out-Null
$hls3.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d4040b3e1f70ff1b3206d7be231a7cda26b7619b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a8aa8e86-e13d-403a-8a83-e5f5fb1a2894.01c7048035d0e34fda3ee283b56e83b8ec8c6c1a.de-de.xlf", "", "", "a8aa8e86-e13d-403a-8a83-e5f5fb1a2894.01c7048035d0e34fda3ee283b56e83b8ec8c6c1a.de-de.xlf") | Out-Null
$hls3.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9c7cda88b83d9c32733cab49e364485c4000bf43/e2e/99319a3b-4a97-4338-b29b-ed4a462b98fc.md", "", "", "99319a3b-4a97-4338-b29b-ed4a462b98fc.md") | Out-Null
$hls3.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac8ecf3ddea2df7b108195a6eaa0a2e43667b30b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/99319a3b-4a97-4338-b29b-ed4a462b98fc.91f6bf9d511a154dc72d665b05a013e920cdbb30.de-de.xlf", "", "", "99319a3b-4a97-4338-b29b-ed4a462b98fc.91f6bf9d511a154dc72d665b05a013e920cdbb30.de-de.xlf") | Out-Null
$hls3.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/9c7cda88b83d9c32733cab49e364485c4000bf43/.localization-config", "", "", ".localization-config") | Out-Null
